# Applies:
#  1) Table style change on the table in slide 16 (graphicFrame shape) from
#     {2F167031-07CB-46E6-9D83-D0AF9EE54417} to {DE1BBE6E-8BF2-4256-983B-8C977B7680A7}.
#  2) Theme color swap: ppt/theme/theme1.xml (the deck's slide-master theme,
#     currently "Integral") receives the "Office Theme" color values.

$p = $ppt.ActivePresentation

# --- 1) Table style ---------------------------------------------------
$s16 = $p.Slides.Item(16)
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $candidate = $s16.Shapes.Item($i)
    if ($candidate.HasTable) {
        $candidate.Table.ApplyStyle("{DE1BBE6E-8BF2-4256-983B-8C977B7680A7}")
    }
}

# --- 2) Theme colors ----------------------------------------------------
# Office Theme color scheme (dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink), in
# RGB() long form (R + G*256 + B*65536) so it lines up with ThemeColor.RGB.
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
